# Update column F (dSF) values on Sheet1 for a subset of rows,
# reflecting a repull/recalculation of the underlying data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -3
    7  = -4
    9  = -3
    13 = -5
    16 = -7
    22 = 6
    23 = 2
    24 = -1
    25 = -4
    28 = -5
    34 = -4
    36 = 3
    40 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
